$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.104.38'
$ws.Range('E2').Value = '  +3.75%  '

# Row 3
$ws.Range('D3').Value = '2.694.03'
$ws.Range('E3').Value = '  +2.16%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.31%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.97%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.55%  '

# Row 7
$ws.Range('E7').Value = '  +0.11%  '

# Row 8
$ws.Range('E8').Value = '  +2.20%  '

# Row 9
$ws.Range('D9').Value = '2.714.59'
$ws.Range('E9').Value = '  +1.89%  '

# Row 10
$ws.Range('E10').Value = '  +4.93%  '

# Row 11
$ws.Range('E11').Value = '  +1.45%  '

# Row 12
$ws.Range('E12').Value = '  +1.20%  '

# Row 13
$ws.Range('E13').Value = '  +2.72%  '

# Row 14
$ws.Range('D14').Value = '3.165.48'
$ws.Range('E14').Value = '  +2.12%  '

# Row 15
$ws.Range('D15').Value = '60.999.53'
$ws.Range('E15').Value = '  +3.57%  '

# Row 16
$ws.Range('D16').Value = '2.938.24'
$ws.Range('E16').Value = '  +10.37%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '21.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.59%  '

# Row 18
$ws.Range('E18').Value = '  +1.22%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '349.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.39%  '

# Row 20
$ws.Range('E20').Value = '  +0.08%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.69%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.82%  '

# Row 23
$ws.Range('E23').Value = '  -0.03%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.37%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.423'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.98%  '

# Row 26
$ws.Range('E26').Value = '  +5.58%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.994'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.94%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.00%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0822'
$ws.Range('E29').Value = '  +2.58%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.17%  '

# Row 31
$ws.Range('E31').Value = '  +0.01%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.43%  '

# Row 33
$ws.Range('E33').Value = '  +2.11%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.26'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.32%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.10%  '

# Row 36
$ws.Range('E36').Value = '  +10.19%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.952'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.00%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.884'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.54%  '

# Row 39
$ws.Range('E39').Value = '  +8.89%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.03%  '

# Row 41
$ws.Range('E41').Value = '  -0.89%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '284.06'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.06%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.14'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.97%  '

# Row 44
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0992'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.70%  '

# Row 45
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.614'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.26%  '

# Row 46
$ws.Range('D46').Value = '2.149.00'
$ws.Range('E46').Value = '  +8.30%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.14%  '

# Row 48
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0541'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.33%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.64%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0236'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.80%  '
